$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# D-column values that are not plain numbers (e.g. "30.606.60") are
# naturally kept as text; values that look like plain numbers need the
# cell formatted as Text first so they stay as text instead of becoming numeric.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.606.60"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.882.84"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.12"
$ws.Range("E5").Value = "  +1.42%  "

$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4756"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2931"
$ws.Range("E8").Value = "  +1.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06532"
$ws.Range("E9").Value = "  +0.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.99"
$ws.Range("E10").Value = "  +1.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07741"
$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7445"
$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.72"
$ws.Range("E13").Value = "  -0.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.883.85"
$ws.Range("E14").Value = "  +0.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.300"
$ws.Range("E15").Value = "  +3.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "275.20"
$ws.Range("E16").Value = "  +0.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.595.54"
$ws.Range("E17").Value = "  +0.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.22"
$ws.Range("E18").Value = "  -3.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007546"
$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.131.07"
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.357"
$ws.Range("E22").Value = "  +1.95%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.234"
$ws.Range("E24").Value = "  +1.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.224"
$ws.Range("E25").Value = "  -0.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.91"
$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.87"
$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.919"
$ws.Range("E28").Value = "  -1.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.343"
$ws.Range("E29").Value = "  -2.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09714"
$ws.Range("E30").Value = "  -2.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.510"
$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.299"
$ws.Range("E32").Value = "  -0.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.164"
$ws.Range("E33").Value = "  +2.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04866"
$ws.Range("E34").Value = "  +2.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.129"
$ws.Range("E35").Value = "  +0.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7001"
$ws.Range("E36").Value = "  +0.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01898"
$ws.Range("E38").Value = "  +1.53%  "

$ws.Range("E39").Value = "  +1.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.329"
$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.92"
$ws.Range("E41").Value = "  +6.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.008"
$ws.Range("E42").Value = "  +4.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4240"
$ws.Range("E43").Value = "  +1.67%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8421"
$ws.Range("E44").Value = "  +0.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.0000"
$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.37"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.381"
$ws.Range("E47").Value = "  +0.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.66"
$ws.Range("E48").Value = "  +0.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.046"
$ws.Range("E49").Value = "  -0.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "918.41"
$ws.Range("E50").Value = "  -1.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05734"
$ws.Range("E51").Value = "  +2.06%  "
